$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# CU "CRU renta de espacio" terminado -> mark the "Eliminar" (K) column
# progress cell for row 8 as complete. The dependent running-total
# formulas across the sheet recalculate automatically.
$ws.Range("K8").Value = 1

# Move the active selection in the frozen bottom-right pane to K9.
$ws.Range("K9").Select()
